$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2267693919272773
$ws.Range("A2").Value = -0.00599999994276601
$ws.Range("A3").Value = -0.054267379449688846
$ws.Range("A4").Value = -0.0079999999032835944
$ws.Range("A5").Value = -0.0029999999434258129
$ws.Range("A6").Value = -0.0019999999378477185
$ws.Range("A7").Value = -0.00999999986268163
$ws.Range("A8").Value = -0.0099999998580160288
$ws.Range("A9").Value = -0.001999999925454965
$ws.Range("A10").Value = -0.001999999919904738
$ws.Range("A11").Value = -0.002999999910118234
$ws.Range("A12").Value = -0.0034999999042524266
$ws.Range("A13").Value = -0.0034999998984917013
$ws.Range("A14").Value = -0.0079999998555964069
$ws.Range("A15").Value = -0.00099999991858901183
$ws.Range("A16").Value = -0.0019999999084143738
$ws.Range("A17").Value = -0.001999999906892036
$ws.Range("A18").Value = -0.0039999998881734555
$ws.Range("A19").Value = -0.0039999999610644821
$ws.Range("A20").Value = -0.0039999999583066881
$ws.Range("A21").Value = -0.0039999999578750334
$ws.Range("A22").Value = -0.0039999999575321965
$ws.Range("A23").Value = -0.0049999999332568379
$ws.Range("A24").Value = -0.01999999979063638
$ws.Range("A25").Value = -0.019999999787942535
$ws.Range("A26").Value = -0.0024999999325672206
$ws.Range("A27").Value = -0.0024999999314014865
$ws.Range("A28").Value = -0.0019999999303959015
$ws.Range("A29").Value = -0.0069999998808532027
$ws.Range("A30").Value = -0.059999999394992454
$ws.Range("A31").Value = 0.00071251223041990386
$ws.Range("A32").Value = -0.0099999998593318651
$ws.Range("A33").Value = 0.047993854734903252
